$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = "B13","C13","D13","B14","C14","D14","C16","D16"
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

# Row 13: Enterprises density (per 1000 people)
$ws.Range("B13").Value = "40.25"
$ws.Range("C13").Value = "0.69"
$ws.Range("D13").Value = "40.94"

# Row 14: Employment (% of total)
$ws.Range("B14").Value = "76.94"
$ws.Range("C14").Value = "12.67"
$ws.Range("D14").Value = "89.61"

# Row 16: Enterprises (% of total)
$ws.Range("C16").Value = "1.67"
$ws.Range("D16").Value = "99.78"
